# Remove the four "show..." command rows (showCurve, showExtraCurve,
# showEvents, showBackgroundEvents) from the "Commands" sheet. These were
# documentation rows for RC button-event display toggles that have been
# dropped from the spec; deleting the rows also removes the now-unused
# shared strings and shifts every following row (and its shared-string
# indices) up by four.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Rows 100-103 (1-based) hold:
#   100: showCurve(<name>,<bool>)
#   101: showExtraCurve(<extra_device>,<curve>,<bool>)
#   102: showEvents(<event_type>, <bool>)
#   103: showBackgroundEvents(<bool>)
$ws.Rows("100:103").Delete()

# Leave the selection where the author ended up after the edit.
$ws.Cells.Item(99, 3).Select() | Out-Null
